$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Create the three new character styles referenced by the diff.
# ---------------------------------------------------------------------------

$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# ---------------------------------------------------------------------------
# 2. Apply "GaNStyle" to every "Datas das campanhas de 2022..." run (4x).
# ---------------------------------------------------------------------------

$campaignText = "Datas das campanhas de 2022 que usam Constelação de Gêmeos: 14 a 23 de fevereiro, 14 a 24 de março"

$rng = $d.Content
$rng.Start = 0
$campaignCount = 0
while ($rng.Find.Execute($campaignText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
    $campaignCount = $campaignCount + 1
    $rng.Collapse(0)
    $rng.End = $d.Content.End
}
Write-Output "GaNStyle applied to $campaignCount run(s)"

# ---------------------------------------------------------------------------
# 3. Apply "GaNParagraph" to the "Está a participar numa campanha global..."
#    paragraph.
# ---------------------------------------------------------------------------

$participarText = "Está a participar numa campanha global para observar e registar as estrelas mais fracas visíveis como forma de medir a poluição luminosa num determinado local. Localizando e observando a  Constelação de Gêmeos no céu noturno e,  comparando-a com cartas estelares, pessoas de todo o mundo aprenderão  como as luzes da sua comunidade contribuem para a poluição luminosa. As suas contribuições para a base de dados on-line irão documentar a visibilidade do céu noturno em todo o mundo."

$rng2 = $d.Content
$rng2.Start = 0
if ($rng2.Find.Execute($participarText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng2.Style = "GaNParagraph"
    Write-Output "GaNParagraph applied"
} else {
    Write-Output "GaNParagraph: text not found!"
}

# ---------------------------------------------------------------------------
# 4. Apply "GaNLinks" to the "por Jenik Hollan, CzechGlobe (...)" run.
# ---------------------------------------------------------------------------

$jenikText = "por Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$rng3 = $d.Content
$rng3.Start = 0
if ($rng3.Find.Execute($jenikText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng3.Style = "GaNLinks"
    Write-Output "GaNLinks applied"
} else {
    Write-Output "GaNLinks: text not found!"
}
